$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment Schedule")

# Insert a new (blank) column before column N (14) - shifts N,O,P -> O,P,Q
$ws.Columns.Item(14).Insert()

# The inserted column should pick up the same width as column M (13),
# which is 11.140625 characters - closest representable width.
$ws.Columns.Item(14).ColumnWidth = 10.3

# Make "Repayment Schedule" the active sheet/tab and set its selection
$ws.Activate()
$ws.Range("U8").Select()
